$wb = $excel.ActiveWorkbook

# --- Sheet "all_tools": update timeout/correlation stats for rows 9-12 ---
$wsAll = $wb.Worksheets.Item("all_tools")

# Row 9
$wsAll.Range("G9").Value = 1143
$wsAll.Range("I9").Value = -0.166102596545867
$wsAll.Range("J9").Value = 0.01727604806480851
$wsAll.Range("K9").Value = -0.2527623213330977
$wsAll.Range("L9").Value = 0.01117475265921138

# Row 10
$wsAll.Range("G10").Value = 859

# Row 11
$wsAll.Range("G11").Value = 859

# Row 12
$wsAll.Range("G12").Value = 859

# --- Sheet "openjml": update timeout/correlation stats for rows 9-12 and a column width ---
$wsJml = $wb.Worksheets.Item("openjml")

# Column J (10th column) width tweak (narrowed by about one character)
$wsJml.Columns.Item(10).ColumnWidth = 18.86

# Row 9
$wsJml.Range("G9").Value = 751
$wsJml.Range("H9").Value = 100
$wsJml.Range("I9").Value = -0.1230941982211943
$wsJml.Range("J9").Value = 0.107046316578524
$wsJml.Range("K9").Value = -0.1671214681732231
$wsJml.Range("L9").Value = 0.09653017580355105

# Row 10
$wsJml.Range("G10").Value = 215
$wsJml.Range("H10").Value = 50
$wsJml.Range("I10").Value = -0.1410673005708742
$wsJml.Range("J10").Value = 0.1918480862059539
$wsJml.Range("K10").Value = -0.1594915550278049
$wsJml.Range("L10").Value = 0.2685745446816231

# Row 11
$wsJml.Range("G11").Value = 215
$wsJml.Range("H11").Value = 50
$wsJml.Range("I11").Value = -0.1104149035826812
$wsJml.Range("J11").Value = 0.285180706372306
$wsJml.Range("K11").Value = -0.1439811999024739
$wsJml.Range("L11").Value = 0.318492693335327

# Row 12
$wsJml.Range("G12").Value = 215
$wsJml.Range("H12").Value = 50
$wsJml.Range("I12").Value = 0.1113864967082373
$wsJml.Range("J12").Value = 0.2778961820094916
$wsJml.Range("K12").Value = 0.1553560043181458
$wsJml.Range("L12").Value = 0.2813439520692285
